$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I5").Value = 54
$ws.Range("L5").Value = 180.83333
$ws.Range("K5").Value = 54
$ws.Range("H5").Value = 123.181816
$ws.Range("N5").Value = -410.83333
$ws.Range("M5").Value = 61
$ws.Range("J5").Value = 180.83333
$ws.Range("N28").Value = -1468
$ws.Range("H28").Value = 644
$ws.Range("L28").Value = 498
$ws.Range("J28").Value = 498
$ws.Range("J57").Value = 15714.286
$ws.Range("N57").Value = -48140.858
$ws.Range("L57").Value = 47142.858
$ws.Range("H57").Value = 22500
$ws.Range("M98").Value = -6169.143
$ws.Range("I98").Value = 7667.143
$ws.Range("H98").Value = 7958.75
$ws.Range("K98").Value = 7667.143
$ws.Range("J107").Value = 1650
$ws.Range("H107").Value = 2162
$ws.Range("N107").Value = -5490
$ws.Range("L107").Value = 1650
$ws.Range("K113").Value = 55559572
$ws.Range("L113").Value = 5610
$ws.Range("I113").Value = 55559572
$ws.Range("M113").Value = -55556318
$ws.Range("H113").Value = 43482624
$ws.Range("J113").Value = 5610
$ws.Range("N113").Value = -12118
$ws.Range("I122").Value = 7667.143
$ws.Range("M122").Value = -20551.429
$ws.Range("H122").Value = 7958.75
$ws.Range("K122").Value = 23001.429
$ws.Range("I132").Value = 2065
$ws.Range("J132").Value = 2443
$ws.Range("L132").Value = 7329
$ws.Range("N132").Value = -12389
$ws.Range("H132").Value = 2159.5
$ws.Range("M132").Value = -3665
$ws.Range("K132").Value = 6195
$ws.Range("I137").Value = 1963.1818
$ws.Range("M137").Value = -3339.5454
$ws.Range("J137").Value = 2866.7778
$ws.Range("H137").Value = 2369.8
$ws.Range("L137").Value = 8600.3334
$ws.Range("N137").Value = -13700.3334
$ws.Range("K137").Value = 5889.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 7552.909
$ws.Range("H32").Value = 7692.7427
$ws.Range("M32").Value = -7265.909
$ws.Range("K32").Value = 7552.909
$ws.Range("I61").Value = 35721412
$ws.Range("K61").Value = 35721412
$ws.Range("M61").Value = -35721200
$ws.Range("H61").Value = 29420280
$ws.Range("M74").Value = -709.7141999999999
$ws.Range("L74").Value = 15000
$ws.Range("N74").Value = -16748
$ws.Range("H74").Value = 1895.721
$ws.Range("I74").Value = 1583.7142
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 1583.7142
$ws.Range("L77").Value = 75000
$ws.Range("K77").Value = 7918.571
$ws.Range("N77").Value = -83736
$ws.Range("I77").Value = 1583.7142
$ws.Range("H77").Value = 1895.721
$ws.Range("J77").Value = 15000
$ws.Range("M77").Value = -3550.571
$ws.Range("M136").Value = -107161686
$ws.Range("K136").Value = 107164236
$ws.Range("I136").Value = 35721412
$ws.Range("H136").Value = 29420280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J132").Value = 100000
$ws.Range("L132").Value = 100000
$ws.Range("N132").Value = -110120
$ws.Range("H132").Value = 100000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L4").Value = 0
$ws.Range("H4").Value = 1000000
$ws.Range("J4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("I132").Value = 2716
$ws.Range("J132").Value = 3877.2856
$ws.Range("L132").Value = 11631.8568
$ws.Range("N132").Value = -16691.8568
$ws.Range("H132").Value = 3341.3076
$ws.Range("M132").Value = -5618
$ws.Range("K132").Value = 8148

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 100000
$ws.Range("M56").Value = -99470
$ws.Range("K56").Value = 100000
$ws.Range("I56").Value = 100000
$ws.Range("I139").Value = 1817.1666
$ws.Range("H139").Value = 1817.1666
$ws.Range("M139").Value = -311.4997999999996
$ws.Range("N139").ClearContents()
$ws.Range("L139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 5451.4998
$ws.Range("J140").Value = 3441.8333
$ws.Range("N140").Value = -20685.4999
$ws.Range("L140").Value = 10325.4999
$ws.Range("H140").Value = 2100.1738
$ws.Range("M141").Value = -25846.999
$ws.Range("I141").Value = 10342.333
$ws.Range("H141").Value = 10342.333
$ws.Range("K141").Value = 31026.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N3").Value = -150482
$ws.Range("H3").Value = 575247.25
$ws.Range("K3").Value = 1000244.5
$ws.Range("M3").Value = -1000128.5
$ws.Range("I3").Value = 1000244.5
$ws.Range("J3").Value = 150250
$ws.Range("L3").Value = 150250
$ws.Range("L32").Value = 17000000
$ws.Range("N32").Value = -17000592
$ws.Range("J32").Value = 17000000
$ws.Range("H32").Value = 17000000
$ws.Range("K68").Value = 20000
$ws.Range("J68").Value = 0
$ws.Range("I68").Value = 20000
$ws.Range("L68").Value = 0
$ws.Range("H68").Value = 20000
$ws.Range("M68").Value = -19189
$ws.Range("N68").ClearContents()
$ws.Range("L71").Value = 0
$ws.Range("K71").Value = 60000
$ws.Range("N71").ClearContents()
$ws.Range("M71").Value = -55944
$ws.Range("J71").Value = 0
$ws.Range("H71").Value = 20000
$ws.Range("I71").Value = 20000
$ws.Range("M107").Value = 1757
$ws.Range("H107").Value = 268
$ws.Range("K107").Value = 163
$ws.Range("I107").Value = 163
$ws.Range("I122").Value = 4900
$ws.Range("J122").Value = 3500
$ws.Range("M122").Value = -12250
$ws.Range("N122").Value = -15400
$ws.Range("L122").Value = 10500
$ws.Range("H122").Value = 4500
$ws.Range("K122").Value = 14700
$ws.Range("I132").Value = 3076.6829
$ws.Range("H132").Value = 3541.4849
$ws.Range("M132").Value = -6700.048699999999
$ws.Range("K132").Value = 9230.048699999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J16").Value = 6998
$ws.Range("H16").Value = 7332.6665
$ws.Range("L16").Value = 6998
$ws.Range("K16").Value = 7500
$ws.Range("M16").Value = -7330
$ws.Range("N16").Value = -7338
$ws.Range("I16").Value = 7500
$ws.Range("K22").Value = 2033
$ws.Range("H22").Value = 2587.25
$ws.Range("M22").Value = -1738
$ws.Range("I22").Value = 2033
$ws.Range("H27").Value = 2587.25
$ws.Range("K27").Value = 2033
$ws.Range("I27").Value = 2033
$ws.Range("M27").Value = -1926
$ws.Range("J61").Value = 6499.6665
$ws.Range("N61").Value = -6903.6665
$ws.Range("L61").Value = 6499.6665
$ws.Range("I61").Value = 256501
$ws.Range("K61").Value = 256501
$ws.Range("M61").Value = -256299
$ws.Range("H61").Value = 149357.58
$ws.Range("H109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("L109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K113").Value = 256501
$ws.Range("L113").Value = 6499.6665
$ws.Range("I113").Value = 256501
$ws.Range("M113").Value = -254331
$ws.Range("H113").Value = 149357.58
$ws.Range("J113").Value = 6499.6665
$ws.Range("N113").Value = -10839.6665
$ws.Range("I132").Value = 15970.5
$ws.Range("J132").Value = 9560.25
$ws.Range("L132").Value = 28680.75
$ws.Range("N132").Value = -33740.75
$ws.Range("H132").Value = 14139
$ws.Range("M132").Value = -45381.5
$ws.Range("K132").Value = 47911.5
$ws.Range("J137").Value = 85100
$ws.Range("H137").Value = 85100
$ws.Range("L137").Value = 85100
$ws.Range("N137").Value = -95300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 19585.416
$ws.Range("J31").Value = 19585.416
$ws.Range("N31").Value = -20281.416
$ws.Range("L31").Value = 19585.416
$ws.Range("H70").Value = 39208.89
$ws.Range("L70").Value = 40412.855
$ws.Range("N70").Value = -41042.855
$ws.Range("J70").Value = 40412.855
$ws.Range("H73").Value = 39208.89
$ws.Range("N73").Value = -42596.855
$ws.Range("L73").Value = 40412.855
$ws.Range("J73").Value = 40412.855
$ws.Range("H106").Value = 41246.375
$ws.Range("L106").Value = 41246.375
$ws.Range("N106").Value = -43770.375
$ws.Range("J106").Value = 41246.375
$ws.Range("I132").Value = 3903.2778
$ws.Range("J132").Value = 4767.3335
$ws.Range("L132").Value = 14302.0005
$ws.Range("N132").Value = -19362.0005
$ws.Range("H132").Value = 4119.2915
$ws.Range("M132").Value = -9179.8334
$ws.Range("K132").Value = 11709.8334
